$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-20 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-21 Thursday", 2) | Out-Null
$d.Content.Find.Execute("447×3=1341", $true, $false, $false, $false, $false, $true, 1, $false, "374×9=3366", 2) | Out-Null
$d.Content.Find.Execute("920×8=7360", $true, $false, $false, $false, $false, $true, 1, $false, "190×7=1330", 2) | Out-Null
$d.Content.Find.Execute("917×2=1834", $true, $false, $false, $false, $false, $true, 1, $false, "681×9=6129", 2) | Out-Null
$d.Content.Find.Execute("270×8=2160", $true, $false, $false, $false, $false, $true, 1, $false, "260×5=1300", 2) | Out-Null
$d.Content.Find.Execute("788×5=3940", $true, $false, $false, $false, $false, $true, 1, $false, "775×6=4650", 2) | Out-Null
$d.Content.Find.Execute("559×5=2795", $true, $false, $false, $false, $false, $true, 1, $false, "224×4=896", 2) | Out-Null
$d.Content.Find.Execute("323×8=2584", $true, $false, $false, $false, $false, $true, 1, $false, "103×7=721", 2) | Out-Null
$d.Content.Find.Execute("741×8=5928", $true, $false, $false, $false, $false, $true, 1, $false, "610×2=1220", 2) | Out-Null
$d.Content.Find.Execute("230×6=1380", $true, $false, $false, $false, $false, $true, 1, $false, "459×5=2295", 2) | Out-Null
$d.Content.Find.Execute("249×4=996", $true, $false, $false, $false, $false, $true, 1, $false, "697×8=5576", 2) | Out-Null
$d.Content.Find.Execute("985×5=4925", $true, $false, $false, $false, $false, $true, 1, $false, "437×7=3059", 2) | Out-Null
$d.Content.Find.Execute("737×4=2948", $true, $false, $false, $false, $false, $true, 1, $false, "445×9=4005", 2) | Out-Null
$d.Content.Find.Execute("381×5=1905", $true, $false, $false, $false, $false, $true, 1, $false, "264×7=1848", 2) | Out-Null
$d.Content.Find.Execute("889×6=5334", $true, $false, $false, $false, $false, $true, 1, $false, "726×9=6534", 2) | Out-Null
$d.Content.Find.Execute("272×6=1632", $true, $false, $false, $false, $false, $true, 1, $false, "550×5=2750", 2) | Out-Null
$d.Content.Find.Execute("790×7=5530", $true, $false, $false, $false, $false, $true, 1, $false, "659×2=1318", 2) | Out-Null
$d.Content.Find.Execute("319×4=1276", $true, $false, $false, $false, $false, $true, 1, $false, "844×6=5064", 2) | Out-Null
$d.Content.Find.Execute("595×7=4165", $true, $false, $false, $false, $false, $true, 1, $false, "280×6=1680", 2) | Out-Null
$d.Content.Find.Execute("494×9=4446", $true, $false, $false, $false, $false, $true, 1, $false, "766×3=2298", 2) | Out-Null
$d.Content.Find.Execute("529×8=4232", $true, $false, $false, $false, $false, $true, 1, $false, "579×5=2895", 2) | Out-Null
$d.Content.Find.Execute("799×7=5593", $true, $false, $false, $false, $false, $true, 1, $false, "689×3=2067", 2) | Out-Null
$d.Content.Find.Execute("504×4=2016", $true, $false, $false, $false, $false, $true, 1, $false, "326×3=978", 2) | Out-Null
$d.Content.Find.Execute("455×5=2275", $true, $false, $false, $false, $false, $true, 1, $false, "545×8=4360", 2) | Out-Null
$d.Content.Find.Execute("315×3=945", $true, $false, $false, $false, $false, $true, 1, $false, "200×3=600", 2) | Out-Null
$d.Content.Find.Execute("397×2=794", $true, $false, $false, $false, $false, $true, 1, $false, "950×8=7600", 2) | Out-Null
